$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# Values are written with a leading apostrophe so Excel stores them as literal
# text (matching the sheet's existing text-formatted cells) instead of
# auto-converting numeric-looking strings into Number/Percentage cells, then the
# style is reset to Normal so no stray "text-stored-as-number" formatting sticks.

$ws.Range("D2").Value = "'317.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.96%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'47.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.49%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.241"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.61%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07970"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.84%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.589"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.95%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.386"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'38.16%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.638"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.85%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'3.30%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1935"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.60%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09274"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.00%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04577"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'9.67%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1042"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.32%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001315"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.87%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04170"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.35%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005832"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.17%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.331"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.12%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.435"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.13%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3466"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'3.61%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.125"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.60%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1392"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.42%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3098"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.97%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001316"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'3.07%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004242"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-5.40%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001353"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.32%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003541"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02639"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'7.78%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05682"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'7.67%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01082"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'81.22%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008167"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'6.47%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1433"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.37%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007680"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.44%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008470"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'13.93%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3470"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'15.09%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006913"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'4.08%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.17%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05492"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'23.80%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.72%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.17%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.17%"
$ws.Range("E51").Style = "Normal"
